$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.36
$ws.Range("E2").Value = 1.24

# Row 3
$ws.Range("B3").Value = 1.47
$ws.Range("D3").Value = 1.39
$ws.Range("E3").Value = 1.27

# Row 4
$ws.Range("C4").Value = 1.46
$ws.Range("E4").Value = 1.24
$ws.Range("F4").Value = 1.09

# Row 5
$ws.Range("B5").Value = 1.48
$ws.Range("C5").Value = 1.37
$ws.Range("D5").Value = 1.34
$ws.Range("F5").Value = 1.03
$ws.Range("G5").Value = 0.72

# Row 6
$ws.Range("D6").Value = 1.53
$ws.Range("E6").Value = 1.33
$ws.Range("G6").Value = 1.06

# Row 7
$ws.Range("E7").Value = 1.92
$ws.Range("F7").Value = 1.44
